$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.934.25"
$ws.Range("E2").Value = "  -2.97%  "

$ws.Range("D3").Value = "1.881.26"
$ws.Range("E3").Value = "  -3.63%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.60%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.63%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4598"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4067"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.40%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.83"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07968"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9887"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.98%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.43%  "

$ws.Range("D13").Value = "1.901.25"
$ws.Range("E13").Value = "  -1.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.895"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.057"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001028"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06559"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.75%  "

$ws.Range("E20").Value = "  -4.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.57%  "

$ws.Range("D22").Value = "28.893.21"
$ws.Range("E22").Value = "  -2.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.404"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.27%  "

$ws.Range("E24").Value = "  +1.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.200"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.72%  "

$ws.Range("D26").Value = "2.087.49"
$ws.Range("E26").Value = "  -3.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.91%  "

$ws.Range("E28").Value = "  -3.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.074"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.462"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.021"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09322"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.50%  "

$ws.Range("E34").Value = "  -6.36%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.518"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.269"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06032"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02221"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.276"
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.170"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.65%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5771"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1822"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.77%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.255"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07504"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.36%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.40%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.241"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5433"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.896"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.47%  "
